# Weekly driver report update for 2025-04-19
# Pace Univ driver summary - "Bad Drivers" table loses one entry (now 3 rows)
# and the "Good Drivers" table is refreshed with this week's samples (now 9 rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A gets narrower this week (59 -> 45 characters) ---
# (COM ColumnWidth is in characters; the engine rounds through a pixel
#  conversion, so 44.142857 is the input that lands exactly on 45 stored.)
$ws.Columns.Item(1).ColumnWidth = 44.142857

# =====================================================================
# "Bad Drivers" table (rows 1-7): 4 data rows -> 3 data rows.
# Delete the old #3 entry (Qualcomm .967) so everything below shifts up
# by one row, carrying the "Totals:" row's formatting from row 7 to 6.
# =====================================================================
$ws.Rows.Item(3).Delete()

$ws.Range("A3").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.190.0.4"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 522
$ws.Range("D3").Value = 94.90000000000001

$ws.Range("A4").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 23.20.1.1"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 931
$ws.Range("D4").Value = 97.3

$ws.Range("A5").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.60.0.10"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 48
$ws.Range("D5").Value = 98.90000000000001

# Row 6 is now the "Totals:" row (shifted up from 7); refresh the sums
# and drop the now-unused D (no "Good Roaming %" on a totals row).
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = 1501
$ws.Range("D6").ClearContents()

# =====================================================================
# "Good Drivers" table: 11 data rows -> 9 data rows (two Qualcomm
# entries roll off entirely). Remove two rows from the data block so
# the heading/header land on 12/13 and data on 14-22.
# =====================================================================
$ws.Rows.Item(23).Delete()
$ws.Rows.Item(23).Delete()

$ws.Range("A14").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.3.1"
$ws.Range("B14").Value = 10661
$ws.Range("D14").Value = 100
$ws.Range("E14").ClearContents()

$ws.Range("A15").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B15").Value = 56018
$ws.Range("D15").Value = 100
$ws.Range("E15").ClearContents()

$ws.Range("A16").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B16").Value = 34244
$ws.Range("D16").Value = 100
$ws.Range("E16").ClearContents()

$ws.Range("A17").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B17").Value = 442178
$ws.Range("D17").Value = 99.90000000000001
$ws.Range("E17").Value = "2024-11-10"

$ws.Range("A18").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.0.3"
$ws.Range("B18").Value = 14239
$ws.Range("D18").Value = 100
$ws.Range("E18").Value = "2022-05-23"

$ws.Range("A19").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.100.1.1"
$ws.Range("B19").Value = 265400
$ws.Range("D19").Value = 99.90000000000001
$ws.Range("E19").Value = "2022-05-01"

$ws.Range("A20").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B20").Value = 77849
$ws.Range("D20").Value = 99.90000000000001
$ws.Range("E20").Value = "2021-08-18"

$ws.Range("A21").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B21").Value = 59673
$ws.Range("D21").Value = 100
$ws.Range("E21").Value = "2020-08-05"

$ws.Range("A22").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B22").Value = 113652
$ws.Range("D22").Value = 100
$ws.Range("E22").Value = "2019-12-14"

# =====================================================================
# Touch the sheet's bottom-right formatted corner (J27) so the saved
# dimension comes back out as A1:J27, matching the shrunk layout
# (was A1:J30 with the old, taller/wider tables).
# =====================================================================
$ws.Cells.Item(27, 10).Borders.Item(1).LineStyle = 0

Write-Host "Weekly driver report refreshed for 2025-04-19"
